$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.0886444015517949
$ws.Range("H2").Value = -6.415145399695909
$ws.Range("I2").Value = -15.2366315985584

$ws.Range("G3").Value = 0.06874065249478947
$ws.Range("H3").Value = 4.004523737037561

$ws.Range("G4").Value = 0.003288705913968952
$ws.Range("H4").Value = -64.92147138133743

$ws.Range("G5").Value = 0.008997070971942561
$ws.Range("H5").Value = 167.1934060198587

$ws.Range("G6").Value = -0.2316342529534376
$ws.Range("H6").Value = -4.719153903929869

$ws.Range("G7").Value = -0.2262178667003966
$ws.Range("H7").Value = 9.467838671133507

$ws.Range("G8").Value = -0.3641012972977667
$ws.Range("H8").Value = 1.662248448915083

$ws.Range("G9").Value = -0.4015886391308283
$ws.Range("H9").Value = -0.7410975280407197

$ws.Range("G10").Value = -0.05770561978316782
$ws.Range("H10").Value = -456.1000781456293

$ws.Range("G11").Value = 0.1140495497016929
$ws.Range("H11").Value = 809.8872868729192

$ws.Range("G12").Value = 0.2143321659924274
$ws.Range("H12").Value = -5.659009165943698

$ws.Range("G13").Value = 0.2808933758141552
$ws.Range("H13").Value = 6.665433221394814

$ws.Range("G14").Value = -0.04638896380487877
$ws.Range("H14").Value = -384.8713605042979

$ws.Range("G15").Value = 0.01633438544050347
$ws.Range("H15").Value = -19.08349689553528

$ws.Range("G16").Value = 0.1438630100107573
$ws.Range("H16").Value = 21.91872943683123

$ws.Range("G17").Value = 0.1807597067886283
$ws.Range("H17").Value = -17.39830333349245

$ws.Range("G18").Value = 0.05838078500473245
$ws.Range("H18").Value = -3.449538017256835

$ws.Range("G19").Value = 0.067394658579482
$ws.Range("H19").Value = -25.19127631536024

$ws.Range("G20").Value = -0.1527120297921342
$ws.Range("H20").Value = -4.934382427745987

$ws.Range("G21").Value = -0.1709671241885869
$ws.Range("H21").Value = 14.44366958839294

$ws.Range("G22").Value = 0.07051987772539384
$ws.Range("H22").Value = 29.66015149805618

$ws.Range("G23").Value = 0.02956814685209262
$ws.Range("H23").Value = -27.60028958475407

$ws.Range("G24").Value = 0.1301955247480747
$ws.Range("H24").Value = 12.49517111100397

$ws.Range("G25").Value = 0.1343131301266305
$ws.Range("H25").Value = -11.68027506615141

$ws.Range("G26").Value = -0.00479430540158214
$ws.Range("H26").Value = -109.0669378628279

$ws.Range("G27").Value = 0.03343570931646109
$ws.Range("H27").Value = -33.74712809402865

$ws.Range("G28").Value = 0.1470176520617906
$ws.Range("H28").Value = -3.854148886645259

$ws.Range("G29").Value = 0.2177131728358107
$ws.Range("H29").Value = 27.53909461745133

$ws.Range("G30").Value = 0.009390998057026095
$ws.Range("H30").Value = -52.00369845701253

$ws.Range("G31").Value = 0.05560221833922227
$ws.Range("H31").Value = 472.9251598820929

$ws.Range("G32").Value = 0.002485560239313358
$ws.Range("H32").Value = -93.33509749574468

$ws.Range("G33").Value = 0.01640260904322053
$ws.Range("H33").Value = -37.16351530895351

$ws.Range("G34").Value = 0.1085003627584462
$ws.Range("H34").Value = -15.21290053596608

$ws.Range("G35").Value = 0.1342981806298475
$ws.Range("H35").Value = 4.381751586847275

$ws.Range("G36").Value = -0.06301702730100128
$ws.Range("H36").Value = -519.2223195706603

$ws.Range("G37").Value = -0.06813582082890965
$ws.Range("H37").Value = -544.9108519775704

$ws.Range("G38").Value = -0.043818990756388
$ws.Range("H38").Value = -2045.451391819906

$ws.Range("G39").Value = -0.0006388819391132361
$ws.Range("H39").Value = 98.0877072721449

$ws.Range("G40").Value = 0.1540609275574361
$ws.Range("H40").Value = 4.413405605206536

$ws.Range("G41").Value = 0.1261357498004187
$ws.Range("H41").Value = -21.84859340847964

$ws.Range("G42").Value = 0.0233380487336705
$ws.Range("H42").Value = -63.85335187320862

$ws.Range("G43").Value = 0.06301167152618659
$ws.Range("H43").Value = 81.27365683151612

$ws.Range("G44").Value = 0.02617031558141085
$ws.Range("H44").Value = 85.43696303300688

$ws.Range("G45").Value = -0.01949309551712053
$ws.Range("H45").Value = -147.476928360069

$ws.Range("G46").Value = -0.05864209085215774
$ws.Range("H46").Value = 10.9051777425079

$ws.Range("G47").Value = -0.02481344624925357
$ws.Range("H47").Value = 39.93312814886443

$ws.Range("G48").Value = -0.1754285798906397
$ws.Range("H48").Value = -39.2562093620714

$ws.Range("G49").Value = -0.1097737721924263
$ws.Range("H49").Value = 44.41311002885247

$ws.Range("G50").Value = 0.1095416054381532
$ws.Range("H50").Value = 0.6135163804743632

$ws.Range("G51").Value = 0.1544880758584554
$ws.Range("H51").Value = 54.07080845757724

$ws.Range("G52").Value = 0.05169947435970361
$ws.Range("H52").Value = -13.28363748651784

$ws.Range("G53").Value = 0.08215057575192232
$ws.Range("H53").Value = 21.62765115549287

$ws.Range("G54").Value = -0.09712432492805938
$ws.Range("H54").Value = -38.90663445571267

$ws.Range("G55").Value = -0.01329752380200852
$ws.Range("H55").Value = 82.7833019244336

$ws.Range("G56").Value = 0.06330801553388377
$ws.Range("H56").Value = 38.14846127320192

$ws.Range("G57").Value = 0.09896715712799001
$ws.Range("H57").Value = 1814.155500560061
